$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated map coordinates (columns J = x, K = y) for a set of locations.
# Values below reflect locations being moved around on the map.
$updates = @(
    @{ Row = 4;  J = 256;  K = 256 },
    @{ Row = 9;  J = 480;  K = 288 },
    @{ Row = 10; J = 336;  K = 688 },
    @{ Row = 11; J = 1216; K = 768 },
    @{ Row = 12; J = 1264; K = 448 },
    @{ Row = 13; J = 800;  K = 320 },
    @{ Row = 14; J = 208;  K = 624 },
    @{ Row = 15; J = 592;  K = 880 },
    @{ Row = 16; J = 960;  K = 384 },
    @{ Row = 17; J = 192;  K = 1216 },
    @{ Row = 35; J = 1152; K = 848 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 10).Value = $u.J
    $ws.Cells.Item($u.Row, 11).Value = $u.K
}

# Columns J (10) and K (11) now contain wider values (up to 4 digits),
# so their best-fit width grows to match column L's width.
$ws.Columns.Item(10).ColumnWidth = 5
$ws.Columns.Item(11).ColumnWidth = 5
